$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 167
$ws.Cells.Item(2, 3).Value = '{''C'': 1.5, ''solver'': ''saga'', ''max_iter'': 100, ''tol'': 0.001}'
$ws.Cells.Item(2, 4).Value = 0.92830000000000001
$ws.Cells.Item(2, 5).Value = 0.97868232000000011
$ws.Cells.Item(2, 7).Value = 0.92862416000000003

# Row 3
$ws.Cells.Item(3, 1).Value = 132
$ws.Cells.Item(3, 3).Value = '{''C'': 1.0, ''solver'': ''saga'', ''max_iter'': 100, ''tol'': 0.0001}'
$ws.Cells.Item(3, 4).Value = 0.9284
$ws.Cells.Item(3, 5).Value = 0.97864359999999984
$ws.Cells.Item(3, 7).Value = 0.92855755999999989

# Row 6
$ws.Cells.Item(6, 1).Value = 150
$ws.Cells.Item(6, 3).Value = '{''C'': 1.2, ''solver'': ''saga'', ''max_iter'': 100, ''tol'': 0.0001}'
$ws.Cells.Item(6, 4).Value = 0.92849999999999999
$ws.Cells.Item(6, 5).Value = 0.97866092000000005
$ws.Cells.Item(6, 7).Value = 0.92856551999999981

# Row 7
$ws.Cells.Item(7, 1).Value = 114
$ws.Cells.Item(7, 3).Value = '{''C'': 0.5, ''solver'': ''saga'', ''max_iter'': 100, ''tol'': 0.0001}'
$ws.Cells.Item(7, 4).Value = 0.92820000000000003
$ws.Cells.Item(7, 5).Value = 0.97852900000000009
$ws.Cells.Item(7, 7).Value = 0.92856236000000003

# Row 63
$ws.Cells.Item(63, 1).Value = 100
$ws.Cells.Item(63, 3).Value = '{''C'': 0.5, ''solver'': ''liblinear'', ''max_iter'': 1000, ''tol'': 0.0001}'
$ws.Cells.Item(63, 4).Value = 0.97219999999999995
$ws.Cells.Item(63, 5).Value = 0.99590708000000006
$ws.Cells.Item(63, 7).Value = 0.91436088000000004

# Row 64
$ws.Cells.Item(64, 1).Value = 102
$ws.Cells.Item(64, 3).Value = '{''C'': 0.5, ''solver'': ''liblinear'', ''max_iter'': 100, ''tol'': 0.0001}'
$ws.Cells.Item(64, 4).Value = 0.97219999999999995
$ws.Cells.Item(64, 5).Value = 0.99590708000000006
$ws.Cells.Item(64, 7).Value = 0.91436088000000004

# Row 65
$ws.Cells.Item(65, 1).Value = 104
$ws.Cells.Item(65, 3).Value = '{''C'': 0.5, ''solver'': ''liblinear'', ''max_iter'': 10000, ''tol'': 0.0001}'
$ws.Cells.Item(65, 4).Value = 0.97219999999999995
$ws.Cells.Item(65, 5).Value = 0.99590708000000006
$ws.Cells.Item(65, 7).Value = 0.91436088000000004

# Row 66
$ws.Cells.Item(66, 1).Value = 99
$ws.Cells.Item(66, 3).Value = '{''C'': 0.5, ''solver'': ''liblinear'', ''max_iter'': 1000, ''tol'': 0.001}'
$ws.Cells.Item(66, 4).Value = 0.97230000000000005
$ws.Cells.Item(66, 5).Value = 0.99590836000000005
$ws.Cells.Item(66, 7).Value = 0.91435268000000003

# Row 67
$ws.Cells.Item(67, 1).Value = 101
$ws.Cells.Item(67, 3).Value = '{''C'': 0.5, ''solver'': ''liblinear'', ''max_iter'': 100, ''tol'': 0.001}'
$ws.Cells.Item(67, 4).Value = 0.97230000000000005
$ws.Cells.Item(67, 5).Value = 0.99590836000000005
$ws.Cells.Item(67, 7).Value = 0.91435268000000003

# Row 68
$ws.Cells.Item(68, 1).Value = 103
$ws.Cells.Item(68, 3).Value = '{''C'': 0.5, ''solver'': ''liblinear'', ''max_iter'': 10000, ''tol'': 0.001}'
$ws.Cells.Item(68, 4).Value = 0.97230000000000005
$ws.Cells.Item(68, 5).Value = 0.99590836000000005
$ws.Cells.Item(68, 7).Value = 0.91435268000000003

# Row 73
$ws.Cells.Item(73, 1).Value = 8
$ws.Cells.Item(73, 3).Value = '{''C'': 0.15, ''max_iter'': 10000, ''loss'': ''hinge'', ''tol'': 0.001}'
$ws.Cells.Item(73, 4).Value = 0.96860000000000002
$ws.Cells.Item(73, 5).Value = 0.99161688000000003
$ws.Cells.Item(73, 7).Value = 0.91032539999999995

# Row 74
$ws.Cells.Item(74, 1).Value = 0
$ws.Cells.Item(74, 3).Value = '{''C'': 0.15, ''max_iter'': 1000, ''loss'': ''hinge'', ''tol'': 0.001}'
$ws.Cells.Item(74, 4).Value = 0.96860000000000002
$ws.Cells.Item(74, 5).Value = 0.99161704000000006
$ws.Cells.Item(74, 7).Value = 0.91032335999999991

# Row 75
$ws.Cells.Item(75, 1).Value = 9
$ws.Cells.Item(75, 3).Value = '{''C'': 0.15, ''max_iter'': 10000, ''loss'': ''hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(75, 4).Value = 0.96860000000000002
$ws.Cells.Item(75, 5).Value = 0.99161692000000001
$ws.Cells.Item(75, 7).Value = 0.91032320000000011

# Row 76
$ws.Cells.Item(76, 1).Value = 1
$ws.Cells.Item(76, 3).Value = '{''C'': 0.15, ''max_iter'': 1000, ''loss'': ''hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(76, 4).Value = 0.96860000000000002
$ws.Cells.Item(76, 5).Value = 0.99161692000000001
$ws.Cells.Item(76, 7).Value = 0.91032307999999995

# Row 92
$ws.Cells.Item(92, 1).Value = 21
$ws.Cells.Item(92, 3).Value = '{''C'': 0.2, ''max_iter'': 10000, ''loss'': ''hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(92, 4).Value = 0.97140000000000004
$ws.Cells.Item(92, 5).Value = 0.99244843999999999
$ws.Cells.Item(92, 7).Value = 0.90721028000000004

# Row 93
$ws.Cells.Item(93, 1).Value = 20
$ws.Cells.Item(93, 3).Value = '{''C'': 0.2, ''max_iter'': 10000, ''loss'': ''hinge'', ''tol'': 0.001}'
$ws.Cells.Item(93, 4).Value = 0.97150000000000003
$ws.Cells.Item(93, 5).Value = 0.99244887999999998
$ws.Cells.Item(93, 7).Value = 0.90721007999999992

# Row 94
$ws.Cells.Item(94, 1).Value = 13
$ws.Cells.Item(94, 3).Value = '{''C'': 0.2, ''max_iter'': 1000, ''loss'': ''hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(94, 4).Value = 0.97140000000000004
$ws.Cells.Item(94, 5).Value = 0.99244880000000013
$ws.Cells.Item(94, 7).Value = 0.90720900000000004

# Row 95
$ws.Cells.Item(95, 1).Value = 12
$ws.Cells.Item(95, 3).Value = '{''C'': 0.2, ''max_iter'': 1000, ''loss'': ''hinge'', ''tol'': 0.001}'
$ws.Cells.Item(95, 4).Value = 0.97140000000000004
$ws.Cells.Item(95, 5).Value = 0.99244868000000008
$ws.Cells.Item(95, 7).Value = 0.90720876000000017

# Row 98
$ws.Cells.Item(98, 1).Value = 136
$ws.Cells.Item(98, 3).Value = '{''C'': 1.2, ''solver'': ''liblinear'', ''max_iter'': 1000, ''tol'': 0.0001}'
$ws.Cells.Item(98, 4).Value = 0.98009999999999997
$ws.Cells.Item(98, 5).Value = 0.9978302
$ws.Cells.Item(98, 7).Value = 0.9064704400000001

# Row 99
$ws.Cells.Item(99, 1).Value = 138
$ws.Cells.Item(99, 3).Value = '{''C'': 1.2, ''solver'': ''liblinear'', ''max_iter'': 100, ''tol'': 0.0001}'
$ws.Cells.Item(99, 4).Value = 0.98009999999999997
$ws.Cells.Item(99, 5).Value = 0.9978302
$ws.Cells.Item(99, 7).Value = 0.9064704400000001

# Row 100
$ws.Cells.Item(100, 1).Value = 140
$ws.Cells.Item(100, 3).Value = '{''C'': 1.2, ''solver'': ''liblinear'', ''max_iter'': 10000, ''tol'': 0.0001}'
$ws.Cells.Item(100, 4).Value = 0.98009999999999997
$ws.Cells.Item(100, 5).Value = 0.9978302
$ws.Cells.Item(100, 7).Value = 0.9064704400000001

# Row 101
$ws.Cells.Item(101, 1).Value = 135
$ws.Cells.Item(101, 3).Value = '{''C'': 1.2, ''solver'': ''liblinear'', ''max_iter'': 1000, ''tol'': 0.001}'
$ws.Cells.Item(101, 4).Value = 0.98009999999999997
$ws.Cells.Item(101, 5).Value = 0.9978312399999999
$ws.Cells.Item(101, 7).Value = 0.90646947999999994

# Row 102
$ws.Cells.Item(102, 1).Value = 137
$ws.Cells.Item(102, 3).Value = '{''C'': 1.2, ''solver'': ''liblinear'', ''max_iter'': 100, ''tol'': 0.001}'
$ws.Cells.Item(102, 4).Value = 0.98009999999999997
$ws.Cells.Item(102, 5).Value = 0.9978312399999999
$ws.Cells.Item(102, 7).Value = 0.90646947999999994

# Row 103
$ws.Cells.Item(103, 1).Value = 139
$ws.Cells.Item(103, 3).Value = '{''C'': 1.2, ''solver'': ''liblinear'', ''max_iter'': 10000, ''tol'': 0.001}'
$ws.Cells.Item(103, 4).Value = 0.98009999999999997
$ws.Cells.Item(103, 5).Value = 0.9978312399999999
$ws.Cells.Item(103, 7).Value = 0.90646947999999994

# Row 117
$ws.Cells.Item(117, 1).Value = 155
$ws.Cells.Item(117, 3).Value = '{''C'': 1.5, ''solver'': ''liblinear'', ''max_iter'': 100, ''tol'': 0.001}'
$ws.Cells.Item(117, 4).Value = 0.98260000000000003
$ws.Cells.Item(117, 5).Value = 0.99818723999999992
$ws.Cells.Item(117, 7).Value = 0.90438607999999998

# Row 118
$ws.Cells.Item(118, 1).Value = 157
$ws.Cells.Item(118, 3).Value = '{''C'': 1.5, ''solver'': ''liblinear'', ''max_iter'': 10000, ''tol'': 0.001}'
$ws.Cells.Item(118, 4).Value = 0.98260000000000003
$ws.Cells.Item(118, 5).Value = 0.99818723999999992
$ws.Cells.Item(118, 7).Value = 0.90438607999999998

# Row 119
$ws.Cells.Item(119, 1).Value = 154
$ws.Cells.Item(119, 3).Value = '{''C'': 1.5, ''solver'': ''liblinear'', ''max_iter'': 1000, ''tol'': 0.0001}'
$ws.Cells.Item(119, 4).Value = 0.98250000000000004
$ws.Cells.Item(119, 5).Value = 0.99818676000000006
$ws.Cells.Item(119, 7).Value = 0.90438151999999994

# Row 120
$ws.Cells.Item(120, 1).Value = 156
$ws.Cells.Item(120, 3).Value = '{''C'': 1.5, ''solver'': ''liblinear'', ''max_iter'': 100, ''tol'': 0.0001}'
$ws.Cells.Item(120, 4).Value = 0.98250000000000004
$ws.Cells.Item(120, 5).Value = 0.99818676000000006
$ws.Cells.Item(120, 7).Value = 0.90438151999999994

# Row 143
$ws.Cells.Item(143, 1).Value = 19
$ws.Cells.Item(143, 3).Value = '{''C'': 0.2, ''max_iter'': 100, ''loss'': ''squared_hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(143, 4).Value = 0.98409999999999997
$ws.Cells.Item(143, 5).Value = 0.99852680000000005
$ws.Cells.Item(143, 7).Value = 0.89689547999999997

# Row 144
$ws.Cells.Item(144, 1).Value = 15
$ws.Cells.Item(144, 3).Value = '{''C'': 0.2, ''max_iter'': 1000, ''loss'': ''squared_hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(144, 4).Value = 0.98409999999999997
$ws.Cells.Item(144, 5).Value = 0.99854080000000001
$ws.Cells.Item(144, 7).Value = 0.89678199999999997

# Row 151
$ws.Cells.Item(151, 1).Value = 37
$ws.Cells.Item(151, 3).Value = '{''C'': 0.5, ''max_iter'': 1000, ''loss'': ''hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(151, 4).Value = 0.98029999999999995
$ws.Cells.Item(151, 5).Value = 0.99481244000000002
$ws.Cells.Item(151, 7).Value = 0.89586387999999983

# Row 152
$ws.Cells.Item(152, 1).Value = 31
$ws.Cells.Item(152, 3).Value = '{''C'': 0.25, ''max_iter'': 100, ''loss'': ''squared_hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(152, 4).Value = 0.98509999999999998
$ws.Cells.Item(152, 5).Value = 0.99875535999999987
$ws.Cells.Item(152, 7).Value = 0.89501131999999994

# Row 156
$ws.Cells.Item(156, 1).Value = 30
$ws.Cells.Item(156, 3).Value = '{''C'': 0.25, ''max_iter'': 100, ''loss'': ''squared_hinge'', ''tol'': 0.001}'
$ws.Cells.Item(156, 4).Value = 0.98480000000000001
$ws.Cells.Item(156, 5).Value = 0.99873379999999989
$ws.Cells.Item(156, 7).Value = 0.89519508000000014

# Row 157
$ws.Cells.Item(157, 1).Value = 34
$ws.Cells.Item(157, 3).Value = '{''C'': 0.25, ''max_iter'': 10000, ''loss'': ''squared_hinge'', ''tol'': 0.001}'
$ws.Cells.Item(157, 4).Value = 0.98480000000000001
$ws.Cells.Item(157, 5).Value = 0.99873379999999989
$ws.Cells.Item(157, 7).Value = 0.89519508000000014

# Row 158
$ws.Cells.Item(158, 1).Value = 27
$ws.Cells.Item(158, 3).Value = '{''C'': 0.25, ''max_iter'': 1000, ''loss'': ''squared_hinge'', ''tol'': 0.0001}'
$ws.Cells.Item(158, 4).Value = 0.98560000000000003
$ws.Cells.Item(158, 5).Value = 0.99881703999999993
$ws.Cells.Item(158, 7).Value = 0.89446732000000007

# Row 186
$ws.Cells.Item(186, 1).Value = 64
$ws.Cells.Item(186, 3).Value = '{''C'': 2.0, ''max_iter'': 100, ''loss'': ''hinge'', ''tol'': 0.001}'
$ws.Cells.Item(186, 4).Value = 0.99099999999999999
$ws.Cells.Item(186, 5).Value = 0.99794247999999997
$ws.Cells.Item(186, 7).Value = 0.88270351999999996

# Row 187
$ws.Cells.Item(187, 1).Value = 62
$ws.Cells.Item(187, 3).Value = '{''C'': 2.0, ''max_iter'': 1000, ''loss'': ''squared_hinge'', ''tol'': 0.001}'
$ws.Cells.Item(187, 4).Value = 0.99760000000000004
$ws.Cells.Item(187, 5).Value = 0.99990064000000001
$ws.Cells.Item(187, 7).Value = 0.87842915999999993
